$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '23.192.20', '  +12.84%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.627.36', '  +10.65%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '''0.9935', '  -1.57%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '''303.88', '  +9.87%  ')
    ,@(6, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '''0.9851', '  +2.72%  ')
    ,@(7, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '''0.3684', '  +3.48%  ')
    ,@(8, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '''0.3422', '  +11.67%  ')
    ,@(9, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '''42.76', '  +8.52%  ')
    ,@(10, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '''1.158', '  +6.54%  ')
    ,@(11, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '''0.07110', '  +7.33%  ')
    ,@(12, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '''0.9906', '  -1.27%  ')
    ,@(13, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '''20.26', '  +12.29%  ')
    ,@(14, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '''5.928', '  +8.97%  ')
    ,@(15, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '''6.661', '  +8.06%  ')
    ,@(16, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.626.34', '  +10.30%  ')
    ,@(17, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '''0.00001083', '  +6.30%  ')
    ,@(18, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '''0.9840', '  +2.55%  ')
    ,@(19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '''0.06681', '  +12.22%  ')
    ,@(20, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '''78.58', '  +14.13%  ')
    ,@(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '''16.15', '  +11.50%  ')
    ,@(22, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '''6.023', '  +10.09%  ')
    ,@(23, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '''11.77', '  +4.70%  ')
    ,@(24, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '23.198.93', '  +12.84%  ')
    ,@(25, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '''2.382', '  +4.65%  ')
    ,@(26, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '''2.616', '  +25.77%  ')
    ,@(27, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '''149.93', '  +3.29%  ')
    ,@(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '''19.40', '  +13.39%  ')
    ,@(29, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '1.802.68', '  +10.30%  ')
    ,@(30, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '''125.05', '  +9.84%  ')
    ,@(31, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '''4.041', '  +4.40%  ')
    ,@(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '''6.080', '  +23.83%  ')
    ,@(33, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '''0.9798', '  +23.52%  ')
    ,@(34, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '''0.08298', '  +4.66%  ')
    ,@(35, 'WEMIXTOKEN', 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix', '''1.674', '  +16.22%  ')
    ,@(36, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '''11.96', '  +16.70%  ')
    ,@(37, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '''8.724', '  +20.02%  ')
    ,@(38, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '''5.233', '  +11.54%  ')
    ,@(39, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '''0.06196', '  +8.12%  ')
    ,@(40, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '''1.260', '  +1.65%  ')
    ,@(41, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '''0.02240', '  +10.20%  ')
    ,@(42, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '''0.2031', '  +9.71%  ')
    ,@(43, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '''0.5943', '  +13.74%  ')
    ,@(44, 'Frax', 'https://coinranking.com/coin/KfWtaeV1W+frax-frax', '''0.9829', '  +2.39%  ')
    ,@(45, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '''3.807', '  +8.51%  ')
    ,@(46, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '''13.09', '  +8.75%  ')
    ,@(47, 'Decentraland', 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana', '''0.5795', '  +12.23%  ')
    ,@(48, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '''126.33', '  +6.56%  ')
    ,@(49, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '''1.982', '  +10.33%  ')
    ,@(50, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '''0.06972', '  +8.54%  ')
    ,@(51, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '''74.49', '  +11.69%  ')
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}
